$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.963.89'
$ws.Range("E2").Value = '  +0.90%  '

$ws.Range("D3").Value = '2.362.22'
$ws.Range("E3").Value = '  +4.23%  '

$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.25%  '

$ws.Range("D5").Value = '''235.24'
$ws.Range("E5").Value = '  +1.72%  '

$ws.Range("D6").Value = '''0.659'
$ws.Range("E6").Value = '  +2.86%  '

$ws.Range("D7").Value = '''73.19'
$ws.Range("E7").Value = '  +13.52%  '

$ws.Range("E8").Value = '  +0.18%  '

$ws.Range("D9").Value = '''0.532'
$ws.Range("E9").Value = '  +21.33%  '

$ws.Range("D10").Value = '''0.0988'
$ws.Range("E10").Value = '  +3.36%  '

$ws.Range("D11").Value = '''28.14'
$ws.Range("E11").Value = '  +6.50%  '

$ws.Range("D12").Value = '2.712.40'
$ws.Range("E12").Value = '  +4.41%  '

$ws.Range("E13").Value = '  +2.22%  '

$ws.Range("D14").Value = '''16.86'
$ws.Range("E14").Value = '  +12.54%  '

$ws.Range("D15").Value = '''6.65'
$ws.Range("E15").Value = '  +9.83%  '

$ws.Range("D16").Value = '''0.884'
$ws.Range("E16").Value = '  +7.38%  '

$ws.Range("D17").Value = '2.368.44'
$ws.Range("E17").Value = '  +4.48%  '

$ws.Range("D18").Value = '43.850.24'
$ws.Range("E18").Value = '  +1.01%  '

$ws.Range("E19").Value = '  +4.17%  '

$ws.Range("D20").Value = '''76.05'
$ws.Range("E20").Value = '  +4.26%  '

$ws.Range("D21").Value = '''6.32'
$ws.Range("E21").Value = '  +3.58%  '

$ws.Range("D22").Value = '''251.22'
$ws.Range("E22").Value = '  +1.58%  '

$ws.Range("E23").Value = '  +0.00%  '

$ws.Range("D24").Value = '''3.78'
$ws.Range("E24").Value = '  -2.23%  '

$ws.Range("D25").Value = '''2.48'
$ws.Range("E25").Value = '  +2.62%  '

$ws.Range("D26").Value = '''10.25'
$ws.Range("E26").Value = '  +5.58%  '

$ws.Range("E27").Value = '  -1.45%  '

$ws.Range("D28").Value = '''22.51'
$ws.Range("E28").Value = '  +3.69%  '

$ws.Range("D29").Value = '''173.15'
$ws.Range("E29").Value = '  -0.27%  '

$ws.Range("D30").Value = '''1.54'
$ws.Range("E30").Value = '  +8.28%  '

$ws.Range("E31").Value = '  +1.54%  '

$ws.Range("E32").Value = '  +4.75%  '

$ws.Range("E33").Value = '  +4.38%  '

$ws.Range("E34").Value = '  +4.65%  '

$ws.Range("D35").Value = '''5.14'
$ws.Range("E35").Value = '  +4.60%  '

$ws.Range("D36").Value = '''3.76'
$ws.Range("E36").Value = '  +4.28%  '

$ws.Range("D37").Value = '''2.43'
$ws.Range("E37").Value = '  +7.49%  '

$ws.Range("D38").Value = '''6.43'
$ws.Range("E38").Value = '  +0.38%  '

$ws.Range("E39").Value = '  +5.58%  '

$ws.Range("D40").Value = '''19.31'
$ws.Range("E40").Value = '  +12.93%  '

$ws.Range("E41").Value = '  +0.02%  '

$ws.Range("D42").Value = '''8.90'
$ws.Range("E42").Value = '  +1.30%  '

$ws.Range("D43").Value = '''1.17'
$ws.Range("E43").Value = '  +9.53%  '

$ws.Range("D44").Value = '''1.21'
$ws.Range("E44").Value = '  +2.85%  '

$ws.Range("D45").Value = '''98.59'
$ws.Range("E45").Value = '  +1.90%  '

$ws.Range("D46").Value = '''0.0963'
$ws.Range("E46").Value = '  +2.34%  '

$ws.Range("E47").Value = '  -0.91%  '

$ws.Range("D48").Value = '''0.181'
$ws.Range("E48").Value = '  +13.46%  '

$ws.Range("D49").Value = '1.441.69'
$ws.Range("E49").Value = '  +0.93%  '

$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").Value = '2.588.05'
$ws.Range("E50").Value = '  +4.49%  '

$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").Value = '''2.29'
$ws.Range("E51").Value = '  +1.60%  '
